# Table4.xlsx - "wrangling nfishers data" update
#
# The "Area of residence" column had a couple of labels cleaned up and two
# previously-blank rows (Los Angeles / San Diego) filled in, splitting what
# used to be lumped together with the Santa Barbara residence group.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Strip the trailing underscores from "Santa Barbara____"
$ws.Range("A6").Value = "Santa Barbara"

# Rows 7 and 8 were blank placeholders; fill them with the residence areas
# that were split out of the old Santa Barbara catch-all row.
$ws.Range("A7").Value = "Los Angeles"
$ws.Range("A8").Value = "San Diego"

# Trim the long trailing underscores from the "Mexican nationals..." label.
# The shorter text now fits on a single line, so shrink the row height to
# match the other single-line rows instead of the old two-line height.
$ws.Range("A10").Value = "Mexican nationals licensed in California"
$ws.Rows.Item(10).RowHeight = 17

# Leave the cursor/selection on A11, matching where editing left off.
[void]$ws.Range("A11").Select()
